$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Main" worksheet right after "Sheet1" and make it active.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws.Name = "Main"

# Enter the data in the same order the shared-string table in the target
# file implies: URL first, then the headers, then the price and the note.
$ws.Range("A2").Value = "https://www.nseindia.com/live_market/dynaContent/live_watch/get_quote/GetQuoteFO.jsp?underlying=SUNPHARMA&instrument=OPTSTK&strike=400.00&type=PE&expiry=25JUL2019"
$ws.Range("A1").Value = "S.No"
$ws.Range("C1").Value = "Price"
$ws.Range("B1").Value = "Date"

# "26.00" must stay textual (it is a shared string, not a number, in the
# target workbook) so format the cell as Text before typing the value.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "26.00"

$ws.Range("A3").Value = "dsadsckank"

# Match the final selection recorded on the new sheet.
$ws.Range("B3").Select()
